$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows right above the "TaskdoVPC" section header (row 22),
# pushing the existing VPC / sensitivity-analysis task block down to rows 25-29.
$ws.Rows("22:24").Insert()

# Row 22 becomes the new "Sensitivity" section header - reuse the formatting
# of an existing section header row (row 17, untouched by the insert above it)
# by copying it in wholesale, then set the header text.
$ws.Range("A17:F17").Copy($ws.Range("A22:F22"))
$ws.Range("B22").Value = "Sensitivity"

# Rows 23 & 24: sensXls / sensSheet fields. Columns A/B already inherited the
# correct formatting (style of row 21 above) from the row-insert; only the
# D:F "filler" cells need to be brought in from another fully-populated
# description row.
$ws.Range("A23").Value = "sensXls"
$ws.Range("B23").Value = "xlsfilefor sensitivity Parameter definition; if it is empty, sheet is in this xlsfile"
$ws.Range("D18:F18").Copy($ws.Range("D23:F23"))

$ws.Range("A24").Value = "sensSheet"
$ws.Range("B24").Value = "xlssheet for sensitivity Parameter definition; if empty first sheet is taken"
$ws.Range("D18:F18").Copy($ws.Range("D24:F24"))

# Match the committed row heights exactly.
$ws.Rows(22).RowHeight = 49.8
$ws.Rows(23).RowHeight = 26.4
$ws.Rows(24).RowHeight = 26.4

# The trailing column-F formatting got dropped off the last three rows of the
# (shifted-down) task block in the original commit - clear it to match.
$ws.Range("F27").Clear()
$ws.Range("F28").Clear()
$ws.Range("F29").Clear()

Write-Host "done"
